$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First statistics block (rows 1-6) ---
$ws.Range("F1").Value = "Statistik"

$ws.Range("F2").Value = "Anzahl von Tests:"
$ws.Range("G2").Value = 20

$ws.Range("F3").Value = "Anzahl bestandenden Tests"
$ws.Range("G3").Value = 18

$ws.Range("F4").Value = "Anzahl von fehlgeschlagene Tests"
$ws.Range("G4").Value = 2

$ws.Range("F5").Value = "Anzahl Fehler"
$ws.Range("G5").Value = 2

$ws.Range("F6").Value = "Beseitigte Fehler"
$ws.Range("G6").Value = "JA"

# --- Second statistics block (rows 7-11), added after correction ---
$ws.Range("F7").Value = "Statistik Nach Korrektur"
$ws.Range("G7").ClearContents()

$ws.Range("F8").Value = "Anzahl von Tests:"
$ws.Range("G8").Value = 20

$ws.Range("F9").Value = "Anzahl von bestandenden Test"
$ws.Range("G9").Value = 20

$ws.Range("F10").Value = "Anzahl von fehlgeschlagene Tests"
$ws.Range("G10").Value = 0

$ws.Range("F11").Value = "Beseitigte Fehler"
$ws.Range("G11").Value = "N/A"

# Update the selected cell in the sheet view to F15
$ws.Range("F15").Select()
